$wb = $excel.ActiveWorkbook

# Width used by this engine's ColumnWidth setter snaps to a 1/6-character grid
# (rounded input + 5/6 offset). Use the closest achievable input for the two
# target widths that appear in the diff.
$colWidthNarrow = 16.333333333333332   # -> stored width ~17.1667 (closest to target 17.2159881591797)
$colWidthWide   = 39.166666666666664   # -> stored width 40 exactly

$readyForHandoff = "Ready for handoff"
$newHandoffDate  = "2016-10-21 01:29:02"
$newZhHandoffDatetime = "2016-10-21 01:28:50"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8108817796a383bc3e95fd74bd230133f5548e33/e2e/e654c024-5e90-48a9-ae0b-868192dc5606.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c2784c859423bde75337f630b45ad4a4b26a7cad/e2e/e654c024-5e90-48a9-ae0b-868192dc5606.md."

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")

# Status text "Handed back: in sync with en-US" -> "Ready for handoff"
# Applies to zh-cn (E) and de-de (F) columns, for both data rows (2 and 3),
# since both rows previously shared the exact same status text/date.
$wsOverview.Range("E2").Value = $readyForHandoff
$wsOverview.Range("F2").Value = $readyForHandoff
$wsOverview.Range("G2").Value = $newHandoffDate

$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = $newHandoffDate

# Column width changes: E and F from ~29.98 to ~17.22
$wsOverview.Range("E1").ColumnWidth = $colWidthNarrow
$wsOverview.Range("F1").ColumnWidth = $colWidthNarrow

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Status column (C) shares the same underlying text as the Overview status cells.
$wsZhCn.Range("C2").Value = $readyForHandoff
$wsZhCn.Range("C3").Value = $readyForHandoff

# Priority "ht" -> "mt" for both rows (shared string used by both)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# Latest Handoff Datetime "2016-10-21 01:26:36" -> "2016-10-21 01:28:50" (shared by both rows)
$wsZhCn.Range("H2").Value = $newZhHandoffDatetime
$wsZhCn.Range("H3").Value = $newZhHandoffDatetime

# Error Detail for e654c024 row (row 3)
$wsZhCn.Range("P3").Value = $errorDetail

# Column width changes: C (Status) from ~29.98 to ~17.22; P (Error Detail) from ~13.75 to 40
$wsZhCn.Range("C1").ColumnWidth = $colWidthNarrow
$wsZhCn.Range("P1").ColumnWidth = $colWidthWide

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status column (C) shares the same underlying text as the Overview status cells.
$wsDeDe.Range("C2").Value = $readyForHandoff
$wsDeDe.Range("C3").Value = $readyForHandoff

# Priority "ht" -> "mt" for both rows (shared string used by both)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# Latest Handoff Datetime here shares the same underlying shared string as the
# Overview status date, so it becomes the new handoff date (2016-10-21 01:29:02),
# not the zh-cn value.
$wsDeDe.Range("H2").Value = $newHandoffDate
$wsDeDe.Range("H3").Value = $newHandoffDate

# Error Detail for e654c024 row (row 3)
$wsDeDe.Range("P3").Value = $errorDetail

# Column width changes: C (Status) from ~29.98 to ~17.22; P (Error Detail) from ~13.75 to 40
$wsDeDe.Range("C1").ColumnWidth = $colWidthNarrow
$wsDeDe.Range("P1").ColumnWidth = $colWidthWide

$wb.Save()
